$d = $word.ActiveDocument

# --- Paragraph: Implications ---
# 1) First run text before "Atheresthes stomias"
$r1 = $d.Content
$r1.Find.Execute(
  ": The distribution and extent of the cold pool directly influence thermal stratification, and overall, changes in surface and bottom temperature influence the spatial distribution of demersal community composition and benthic trophic structure (Mueter and Litzow, 2008; Spencer, 2008; Kotwicki and Lauth, 2013). When the cold pool is small, thermal stratification is weak and subarctic demersal fishes and invertebrates with warm-water affinity (e.g., arrowtooth flounder",
  $true, $false, $false, $false, $false, $true, 1, $false,
  ": The cold pool has a strong influence on the vertical structure of the eastern Bering Sea. Changes in the extent of the cold pool, and consequent changes in surface and bottom temperature, influence the spatial structure of the demersal community and strength of benthic-pelagic coupling (Mueter and Litzow, 2008; Spencer, 2008; Kotwicki and Lauth, 2013). When the cold pool is small, species with warm-water affinity (e.g., arrowtooth flounder",
  2) | Out-Null

# 2) Run text after "Atheresthes stomias" up to "flathead sole"
$r2 = $d.Content
$r2.Find.Execute(
  ") are often more diffusely distributed over the eastern Bering Sea shelf as there is no thermal barrier to their advance from the outer to inner shelf. In contrast, the majority of the subarctic fish and invertebrate community is comprised of species with cool-water affinity (e.g., flathead sole",
  $true, $false, $false, $false, $false, $true, 1, $false,
  ") are distribute more widely over the eastern Bering Sea shelf and expand across the shelf and to the north because there is no thermal barrier to migration. In contrast, the distribution of species with cold water affinity (e.g., Bering flounder",
  2) | Out-Null

# 3) Italic species name run: Hippoglossoides elassodon -> Hippoglossoides robustus
$r3 = $d.Content
$r3.Find.Execute(
  "Hippoglossoides elassodon",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Hippoglossoides robustus",
  2) | Out-Null

# 4) The runs " " + "and northern rock sole" collapse into a single run ", Arctic cod"
#    Remove the separate "and northern rock sole" text and the leading space run that precedes it,
#    replacing the pair with a single new text segment.
$r4 = $d.Content
$r4.Find.Execute(
  " and northern rock sole",
  $true, $false, $false, $false, $false, $true, 1, $false,
  ", Arctic cod",
  2) | Out-Null

# 5) Italic species name run: Lepidopsetta polyxystra -> Boreogadus saida
$r5 = $d.Content
$r5.Find.Execute(
  "Lepidopsetta polyxystra",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Boreogadus saida",
  2) | Out-Null

# 6) Final run of the paragraph
$r6 = $d.Content
$r6.Find.Execute(
  "), and these species often contract in their area occupied and shift in mean distribution to the north or northwest when the cold pool is reduced.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  ") contracts to the north when the cold pool is small.",
  2) | Out-Null

# --- Paragraph: "Although the definition of the cold pool boundary..." ---
$r7 = $d.Content
$r7.Find.Execute(
  "Although the definition of the cold pool boundary is the 2°C isotherm, recent studies indicate that the better predictor of spatial distribution for many fishes and crabs could be the 1°C isotherm (Kotwicki and Lauth, 2013) or the 0°C isotherm for pollock",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "While the cold pool area is defined based on the 2°C isotherm, recent studies suggest that a more ecologically relevant temperature for several subarctic fishes and crabs is the 1°C isotherm (Kotwicki and Lauth, 2013) or the 0°C isotherm for walleye pollock",
  2) | Out-Null

# --- Paragraph continuation: "(Baker 2021; Eisner et al. 2020). Given that waters..." ---
$r8 = $d.Content
$r8.Find.Execute(
  "(Baker 2021; Eisner et al. 2020). Given that waters cooler than 1°C and 0°C were much less extensive than those defined by the 2°C isotherm, it would appear that the cold pool produced very little spatial structure in the benthic thermal habitat of the southeastern Bering Sea in 2021, although cooler bottom temperatures in the northern Bering Sea likely provided some spatial structure in the far north-central area of the eastern Bering Sea shelf.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "(Baker 2021; Eisner et al. 2020). Considering the small extent of bottom temperatures cooler than 0°C and 1°C, it is likely that the bottom temperatures on the eastern Bering Sea shelf did not impose a major thermal barrier to migration for subarctic species in 2021. However, cooler bottom temperatures in the northern Bering Sea (Fig. 3) may have imposed some barrier to migration.",
  2) | Out-Null

# --- Paragraph: "Although the mean surface temperature was slightly closer..." ---
$r9 = $d.Content
$r9.Find.Execute(
  "Although the mean surface temperature was slightly closer to its long-term mean than was mean bottom temperature in the eastern Bering Sea, 2021 conditions represent a continuation of the warm phase of surface temperature that has persisted since 2014.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Although the mean surface temperature was closer to its long-term mean than mean bottom temperature in the eastern Bering Sea in 2019, conditions in 2021 are a continuation of above average surface temperatures that has persisted since 2014.",
  2) | Out-Null
